# summer 24 week 5 updates
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nine")

$ws.Range("C2").Value = 10.04

$ws.Range("B3").Value = 9.960000000000001
$ws.Range("E3").Value = 10.65

$ws.Range("F4").Value = 10.23

$ws.Range("C5").Value = 9.300000000000001
$ws.Range("F5").Value = 10.19
$ws.Range("H5").Value = 8.65

$ws.Range("D6").Value = 9.77
$ws.Range("E6").Value = 9.81
$ws.Range("G6").Value = 10.33
$ws.Range("H6").Value = 10.46
$ws.Range("I6").Value = 8.9

$ws.Range("F7").Value = 9.67
$ws.Range("H7").Value = 9.93
$ws.Range("I7").Value = 7.5

$ws.Range("E8").Value = 11.35
$ws.Range("F8").Value = 9.539999999999999
$ws.Range("G8").Value = 10.07

$ws.Range("F9").Value = 11.1
$ws.Range("G9").Value = 12.5
